$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-11-02 Sunday" "2025-11-03 Monday"

Replace-Text "688÷8=" "529÷8="
Replace-Text "268÷8=" "248÷4="
Replace-Text "307÷4=" "476÷6="
Replace-Text "628÷4=" "787÷5="
Replace-Text "646÷8=" "987÷5="
Replace-Text "786÷2=" "545÷6="
Replace-Text "960÷7=" "119÷3="
Replace-Text "982÷9=" "550÷2="
Replace-Text "918÷4=" "633÷5="
Replace-Text "219÷9=" "995÷5="
Replace-Text "544÷8=" "790÷5="
Replace-Text "272÷9=" "513÷9="
Replace-Text "222÷5=" "936÷5="
Replace-Text "674÷4=" "632÷5="
Replace-Text "425÷5=" "923÷5="
Replace-Text "568÷8=" "441÷5="
Replace-Text "617÷8=" "542÷9="
Replace-Text "854÷9=" "581÷7="
Replace-Text "249÷6=" "868÷4="
Replace-Text "994÷8=" "946÷5="
Replace-Text "560÷7=" "596÷3="
Replace-Text "635÷8=" "613÷4="
Replace-Text "209÷3=" "283÷3="
Replace-Text "520÷6=" "663÷5="
Replace-Text "101÷3=" "833÷8="
